# firered-nettux-new-mons.xlsx -- "up to lt surge"
# Adds newly-discovered route/encounter info to several existing Pokemon rows
# and appends the next batch of species (Wingull .. All Starters) at the
# bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B (location info) for existing rows with new route/area data ---
$ws.Cells.Item(55, 2).Value  = 'Route 4, Route 6'
$ws.Cells.Item(58, 2).Value  = 'Route 6'
$ws.Cells.Item(63, 2).Value  = 'Route 4, Route 11'
$ws.Cells.Item(70, 2).Value  = 'Route 5, Route 6'
$ws.Cells.Item(75, 2).Value  = 'Route 11, DiglettsCave_B1F'
$ws.Cells.Item(79, 2).Value  = 'MtMoon_B2F,DiglettsCave_B1F'
$ws.Cells.Item(81, 2).Value  = 'Route 4, Route 5, Route 11'
$ws.Cells.Item(83, 2).Value  = 'DiglettsCave_B1F'
$ws.Cells.Item(85, 2).Value  = 'Route 11, DiglettsCave_B1F'
$ws.Cells.Item(88, 2).Value  = 'Route 3, MtMoon_B2F, Route 6'
$ws.Cells.Item(93, 2).Value  = 'DiglettsCave_B1F'
$ws.Cells.Item(98, 2).Value  = 'Route 11'
$ws.Cells.Item(110, 2).Value = 'Route 11'
$ws.Cells.Item(120, 2).Value = 'DiglettsCave_B1F'
$ws.Cells.Item(135, 2).Value = 'MtMoon_1F, DiglettsCave_B1F'
$ws.Cells.Item(138, 2).Value = 'Route 2, Route 6'
$ws.Cells.Item(140, 2).Value = 'DiglettsCave_B1F'
$ws.Cells.Item(141, 2).Value = 'DiglettsCave_B1F'
$ws.Cells.Item(144, 2).Value = 'Route 6'
$ws.Cells.Item(154, 2).Value = 'Route 3, Route 6'
$ws.Cells.Item(157, 2).Value = 'Route 6, Route 11'
$ws.Cells.Item(167, 2).Value = 'Route 4, Route 5'
$ws.Cells.Item(170, 2).Value = 'Route 24,Route 5'
$ws.Cells.Item(175, 2).Value = 'Route 5'
$ws.Cells.Item(177, 2).Value = 'Route 11'
$ws.Cells.Item(192, 2).Value = 'Route 11'
$ws.Cells.Item(210, 2).Value = 'Route 5'
$ws.Cells.Item(212, 2).Value = 'Route 5'

# --- Append the next batch of species (species in column A, location in column B where known) ---
$ws.Cells.Item(215, 1).Value = 'Wingull'
$ws.Cells.Item(215, 2).Value = 'Route 25'
$ws.Cells.Item(216, 1).Value = 'Pelipper'
$ws.Cells.Item(217, 1).Value = 'Rotom'
$ws.Cells.Item(218, 1).Value = 'Geodude Alola'
$ws.Cells.Item(218, 2).Value = 'Route 5'
$ws.Cells.Item(219, 1).Value = 'Graveler Alola'
$ws.Cells.Item(219, 2).Value = 'DiglettsCave_B1F'
$ws.Cells.Item(220, 1).Value = 'Golem Alola'
$ws.Cells.Item(221, 1).Value = 'Miltank'
$ws.Cells.Item(221, 2).Value = 'Route 5'
$ws.Cells.Item(222, 1).Value = 'Stantler'
$ws.Cells.Item(223, 1).Value = 'Chinchou'
$ws.Cells.Item(224, 1).Value = 'Lantern'
$ws.Cells.Item(226, 1).Value = 'All Starters'

# --- Column B grew wider text ("Route 4, Route 5, Route 11" etc.) -- widen it to fit ---
$ws.Columns("B").ColumnWidth = 32.14

# --- Leave the view scrolled/selected near the newly-added rows, like the author did ---
$ws.Activate()
$ws.Range("G222").Select() | Out-Null
